$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.269.15"
$ws.Range("E2").Value = "  +3.03%  "
$ws.Range("D3").Value = "1.815.86"
$ws.Range("E3").Value = "  +4.11%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "328.26"
$ws.Range("E5").Value = "  +2.02%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").Value = "0.4372"
$ws.Range("E7").Value = "  +3.49%  "
$ws.Range("D8").Value = "0.3668"
$ws.Range("E8").Value = "  +2.36%  "
$ws.Range("D9").Value = "44.96"
$ws.Range("E9").Value = "  -1.00%  "
$ws.Range("D10").Value = "0.07682"
$ws.Range("E10").Value = "  +3.53%  "
$ws.Range("D11").Value = "1.141"
$ws.Range("E11").Value = "  +2.73%  "
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  +0.03%  "
$ws.Range("E13").Value = "  +3.62%  "
$ws.Range("D14").Value = "6.310"
$ws.Range("E14").Value = "  +3.24%  "
$ws.Range("D15").Value = "7.541"
$ws.Range("E15").Value = "  +4.92%  "
$ws.Range("D16").Value = "1.830.80"
$ws.Range("E16").Value = "  +4.96%  "
$ws.Range("D17").Value = "93.25"
$ws.Range("E17").Value = "  +5.74%  "
$ws.Range("E18").Value = "  +1.67%  "
$ws.Range("D19").Value = "0.06540"
$ws.Range("E19").Value = "  +6.93%  "
$ws.Range("D20").Value = "1.000"
$ws.Range("D21").Value = "17.52"
$ws.Range("E21").Value = "  +3.99%  "
$ws.Range("D22").Value = "6.263"
$ws.Range("E22").Value = "  +2.78%  "
$ws.Range("D23").Value = "28.295.72"
$ws.Range("E23").Value = "  +2.98%  "
$ws.Range("E24").Value = "  +1.70%  "
$ws.Range("D25").Value = "2.043"
$ws.Range("E25").Value = "  -12.63%  "
$ws.Range("D26").Value = "162.46"
$ws.Range("E26").Value = "  +6.25%  "
$ws.Range("D27").Value = "20.74"
$ws.Range("E27").Value = "  +1.84%  "
$ws.Range("D28").Value = "2.036.47"
$ws.Range("E28").Value = "  +4.92%  "
$ws.Range("D29").Value = "2.299"
$ws.Range("E29").Value = "  -3.00%  "
$ws.Range("D30").Value = "128.77"
$ws.Range("E30").Value = "  +2.33%  "
$ws.Range("D31").Value = "1.217"
$ws.Range("E31").Value = "  +2.06%  "
$ws.Range("D32").Value = "5.960"
$ws.Range("E32").Value = "  +5.35%  "
$ws.Range("D33").Value = "0.09199"
$ws.Range("E33").Value = "  +0.87%  "
$ws.Range("D34").Value = "3.488"
$ws.Range("E34").Value = "  -3.83%  "
$ws.Range("D35").Value = "12.94"
$ws.Range("E35").Value = "  +2.56%  "
$ws.Range("D36").Value = "0.02348"
$ws.Range("E36").Value = "  +2.42%  "
$ws.Range("E37").Value = "  +2.20%  "
$ws.Range("D38").Value = "5.202"
$ws.Range("E38").Value = "  +2.55%  "
$ws.Range("D39").Value = "0.6577"
$ws.Range("E39").Value = "  +3.21%  "
$ws.Range("D40").Value = "0.06202"
$ws.Range("E40").Value = "  +2.50%  "
$ws.Range("E41").Value = "  +0.52%  "
$ws.Range("D42").Value = "8.132"
$ws.Range("E42").Value = "  +3.29%  "
$ws.Range("D43").Value = "1.433"
$ws.Range("E43").Value = "  +0.36%  "
$ws.Range("D44").Value = "0.9999"
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("D45").Value = "13.88"
$ws.Range("E45").Value = "  +1.39%  "
$ws.Range("D46").Value = "0.6124"
$ws.Range("E46").Value = "  +4.64%  "
$ws.Range("D47").Value = "3.756"
$ws.Range("E47").Value = "  +1.12%  "
$ws.Range("D48").Value = "125.82"
$ws.Range("E48").Value = "  +0.75%  "
$ws.Range("D49").Value = "2.023"
$ws.Range("E49").Value = "  +4.22%  "
$ws.Range("D50").Value = "1.158"
$ws.Range("E50").Value = "  +3.99%  "
$ws.Range("D51").Value = "0.07006"
$ws.Range("E51").Value = "  +2.38%  "
